$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds dates that must be stored as literal text (shared strings),
# not auto-converted Excel date serials. Typing the text directly (even with
# NumberFormat forced to Text or a leading apostrophe) still changes the
# cell's style in the saved file. Writing the same text via a formula and
# then collapsing it to a static value with Copy/PasteSpecial(xlPasteValues)
# avoids Excel's "looks like a date" auto-conversion while leaving the
# cell's existing style (border formatting) completely untouched.
$ws.Range("A77").Formula = "=""02-12-2025"""
$ws.Range("A77").Copy() | Out-Null
$ws.Range("A77").PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws.Range("B77").Value = "The price of gold in India today is ₹13,020 per gram for 24 karat gold, ₹11,935 per gram for 22 karat gold and ₹9,765 per gram for 18 karat gold (also called 999 gold)."

$excel.CutCopyMode = 0
